# "changes to the observations to add distance travelled"
#
# 1) Zero out the A-matrix / shortest-path observation blocks (rows 4-17,
#    columns C:BN) which previously held 1s / 0.5s.
# 2) Append a new "distance_travelled" observation-modality block: six new
#    rows (18-23) mirroring the layout of the existing "shortest_path"
#    block (row 12-17): label in column A (merged A18:A23), row labels
#    distancet0..distancet20 / distance25t in column B, and zeros across
#    C:BN.
# 3) Update the view (selected cell) to reflect where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the observation matrices back to 0 --------------------------
$ws.Range("C4:BN17").Value = 0

# --- 2. Add the new "distance_travelled" block (rows 18-23) ---------------

# Merge the label column first, while the new rows still carry the sheet's
# default (unstyled) formatting -- merging after the bordered header style
# is pasted in causes the engine to mint extra per-row border variants.
$ws.Range("A18:A23").Merge()

# Row / column labels (these add the required shared strings).
$ws.Range("A18").Value = "distance_travelled"
$ws.Range("B18").Value = "distancet0"
$ws.Range("B19").Value = "distancet5"
$ws.Range("B20").Value = "distancet10"
$ws.Range("B21").Value = "distancet15"
$ws.Range("B22").Value = "distancet20"
$ws.Range("B23").Value = "distance25t"

# Data cells: all zero, same as the other observation-matrix rows.
$ws.Range("C18:BN23").Value = 0

# Copy the formatting (borders / bold header font / wrap / alignment) of the
# existing "shortest_path" block (row 13 is a plain, non-merged row within
# that block) onto the new rows so the new block matches the sheet's style.
$ws.Range("A13:BN13").Copy()
$ws.Range("A18:BN23").PasteSpecial(-4122)

# Match the row height used by the rest of the "shortest_path"-style rows.
$ws.Range("A18:A23").RowHeight = 30

# --- 3. View state ----------------------------------------------------------
$ws.Range("C29").Select()
